# Auto-generated Excel COM-interop script
# Applies scheduled-runner market data refresh updates to Hades_Profits workbook
# Columns: H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#          K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 403
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()

$ws.Range("H137").Value = 3335743.5
$ws.Range("I137").Value = 6668967
$ws.Range("J137").Value = 2520.2666
$ws.Range("K137").Value = 20006901
$ws.Range("L137").Value = 7560.7998
$ws.Range("M137").Value = -20004351
$ws.Range("N137").Value = -12660.7998


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4848261.5
$ws.Range("I32").Value = 5395003.5
$ws.Range("J32").Value = 5688.2856
$ws.Range("K32").Value = 5395003.5
$ws.Range("L32").Value = 5688.2856
$ws.Range("M32").Value = -5394716.5
$ws.Range("N32").Value = -6262.2856

$ws.Range("H35").Value = 3100.5
$ws.Range("I35").Value = 1200
$ws.Range("K35").Value = 1200
$ws.Range("M35").Value = -794

$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("N36").ClearContents()

$ws.Range("H45").Value = 1675.5
$ws.Range("I45").Value = 1543
$ws.Range("J45").Value = 2020
$ws.Range("K45").Value = 1543
$ws.Range("L45").Value = 2020
$ws.Range("M45").Value = -1166
$ws.Range("N45").Value = -2774

$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").ClearContents()


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 1424.6666
$ws.Range("I36").Value = 1424.6666
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 1424.6666
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -890.6666
$ws.Range("N36").ClearContents()

$ws.Range("H135").Value = 60464.4
$ws.Range("J135").Value = 60464.4
$ws.Range("L135").Value = 60464.4
$ws.Range("N135").Value = -70604.39999999999


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 20834428
$ws.Range("I16").Value = 1036.3572
$ws.Range("J16").Value = 166668180
$ws.Range("K16").Value = 1036.3572
$ws.Range("L16").Value = 166668180
$ws.Range("M16").Value = -749.3571999999999
$ws.Range("N16").Value = -166668754

$ws.Range("H31").Value = 8380.692999999999
$ws.Range("I31").Value = 33511.94
$ws.Range("J31").Value = 1014.63794
$ws.Range("K31").Value = 33511.94
$ws.Range("L31").Value = 1014.63794
$ws.Range("M31").Value = -33216.94
$ws.Range("N31").Value = -1604.63794

$ws.Range("H34").Value = 8380.692999999999
$ws.Range("I34").Value = 33511.94
$ws.Range("J34").Value = 1014.63794
$ws.Range("K34").Value = 33511.94
$ws.Range("L34").Value = 1014.63794
$ws.Range("M34").Value = -33309.94
$ws.Range("N34").Value = -1418.63794

$ws.Range("H113").Value = 20834428
$ws.Range("I113").Value = 1036.3572
$ws.Range("J113").Value = 166668180
$ws.Range("K113").Value = 1036.3572
$ws.Range("L113").Value = 166668180
$ws.Range("M113").Value = 1133.6428
$ws.Range("N113").Value = -166672520

$ws.Range("H132").Value = 145450.28
$ws.Range("I132").Value = 3230.4
$ws.Range("J132").Value = 501000
$ws.Range("K132").Value = 9691.200000000001
$ws.Range("L132").Value = 1503000
$ws.Range("M132").Value = -7161.200000000001
$ws.Range("N132").Value = -1508060

$ws.Range("H134").Value = 91769.664
$ws.Range("I134").Value = 1819.4286
$ws.Range("K134").Value = 5458.2858
$ws.Range("M134").Value = -2923.2858

$ws.Range("H141").Value = 72409.62
$ws.Range("J141").Value = 72409.62
$ws.Range("L141").Value = 72409.62
$ws.Range("N141").Value = -82769.62


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 907.7059
$ws.Range("J68").Value = 1000.44446
$ws.Range("L68").Value = 3001.33338
$ws.Range("N68").Value = -4623.33338

$ws.Range("H71").Value = 907.7059
$ws.Range("J71").Value = 1000.44446
$ws.Range("L71").Value = 9004.00014
$ws.Range("N71").Value = -17116.00014

$ws.Range("H75").Value = 3007.5
$ws.Range("J75").Value = 3007.5
$ws.Range("L75").Value = 9022.5
$ws.Range("N75").Value = -11018.5

$ws.Range("H78").Value = 3007.5
$ws.Range("J78").Value = 3007.5
$ws.Range("L78").Value = 27067.5
$ws.Range("N78").Value = -37051.5

$ws.Range("H86").Value = 500
$ws.Range("J86").Value = 500
$ws.Range("L86").Value = 1500
$ws.Range("N86").Value = -3872

$ws.Range("H89").Value = 500
$ws.Range("J89").Value = 500
$ws.Range("L89").Value = 4500
$ws.Range("N89").Value = -16356

$ws.Range("H131").Value = 924.5217
$ws.Range("I131").Value = 586
$ws.Range("J131").Value = 965.8049
$ws.Range("K131").Value = 1758
$ws.Range("L131").Value = 2897.4147
$ws.Range("M131").Value = 3282
$ws.Range("N131").Value = -12977.4147


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2016.7693
$ws.Range("I102").Value = 1989.3334
$ws.Range("K102").Value = 1989.3334
$ws.Range("M102").Value = -367.3334

$ws.Range("H122").Value = 1910.5
$ws.Range("I122").Value = 1499.4
$ws.Range("J122").Value = 2321.6
$ws.Range("K122").Value = 4498.200000000001
$ws.Range("L122").Value = 6964.799999999999
$ws.Range("M122").Value = -2048.200000000001
$ws.Range("N122").Value = -11864.8


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H138").Value = 47800
$ws.Range("J138").Value = 47800
$ws.Range("L138").Value = 47800
$ws.Range("N138").Value = -58080


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 39750
$ws.Range("J46").Value = 39750
$ws.Range("L46").Value = 39750
$ws.Range("N46").Value = -40212

$ws.Range("H107").Value = 1198.8125
$ws.Range("I107").Value = 1298.4445
$ws.Range("J107").Value = 1070.7142
$ws.Range("K107").Value = 3895.3335
$ws.Range("L107").Value = 3212.1426
$ws.Range("M107").Value = -1975.3335
$ws.Range("N107").Value = -7052.142599999999

$ws.Range("H122").Value = 2029.1428
$ws.Range("I122").Value = 1700.6666
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 5101.9998
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -2651.9998
$ws.Range("N122").Value = -16900

$ws.Range("H134").Value = 39750
$ws.Range("J134").Value = 39750
$ws.Range("L134").Value = 119250
$ws.Range("N134").Value = -124320

